$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: DATE_TYPE_CODE - text code, force text format so "001" keeps its
# leading zero instead of being coerced to the number 1.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "001"

# M2 / N2: NOTICE_DATE / REPORT_DATE - stored as plain text strings
# (not Excel date serials) in the source data.
$ws.Range("M2").Value = "2020-12-18 00:00:00"
$ws.Range("N2").Value = "2017-12-31 00:00:00"

# Numeric metric columns updated to the new reported figures.
$ws.Range("O2").Value = -45009866.89
$ws.Range("P2").Value = -90.75140966959999
$ws.Range("Q2").Value = 344153411.22
$ws.Range("R2").Value = 693.9013458348001
$ws.Range("S2").Value = 53772301.73
$ws.Range("T2").Value = 108.418720613
$ws.Range("U2").Value = -139974509.03
$ws.Range("V2").Value = -282.224429664
$ws.Range("Y2").Value = 101327631.52
$ws.Range("Z2").Value = 204.3024348727
$ws.Range("AA2").Value = 135387496.87
$ws.Range("AB2").Value = 272.9758393336
$ws.Range("AC2").Value = -49596879.05

# AD2 (CCE_ADD_RATIO) was empty before; now populated with a number.
$ws.Range("AD2").Value = -158.7750411376
